$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row whose URL (column A) is the "cannes-lion-2019" article. This
# is row 8 in the original sheet; deleting it shifts rows 9-66 up by one.
$ws.Rows.Item(8).Delete()

# After the shift, the "stacey-kennedy" article row (originally row 29) is
# now row 28; clear its "Secondary topic" (column C) value.
$ws.Range("C28").Value = ""

# The "better/the-public-supports-governments..." article row (originally
# row 31) is now row 30; update its "Secondary topic" (column C) value.
$ws.Range("C30").Value = "Leadership content"

# Update the view state to match (pane frozen on row 1, scrolled so row 30
# is the top-left visible cell, with A30 selected).
$ws.Application.ActiveWindow.ScrollRow = 30
$ws.Range("A30").Select()
